$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the explicit (direct) paragraph-mark color override that was
#    stored in <w:pPr><w:rPr><w:color .../></w:rPr></w:pPr> for every
#    Heading 3 / Heading 4 paragraph. The heading color now comes purely
#    from the style definition (see part 2 below), so the direct formatting
#    on the paragraph mark is removed.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $styleName = $p.Style.NameLocal
    if ($styleName -eq "Heading 3" -or $styleName -eq "Heading 4") {
        $full = $p.Range
        $x = $full.XML()
        $mPara = [regex]::Match($x, '<w:p\b.*?</w:p>')
        if ($mPara.Success) {
            $paraXml = $mPara.Value
            $newParaXml = $paraXml -replace '<w:rPr>\s*<w:color\b[^/]*/>\s*</w:rPr>(?=</w:pPr>)', ''
            if ($newParaXml -ne $paraXml) {
                $wrapper = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
                    '<w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
                $full.InsertXML($wrapper)
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Fix the heading styles' colors:
#      Heading 3 / Heading 3 Char : C00000                -> 002060
#      Heading 4 / Heading 4 Char : E36C0A (themed accent6) -> 002060
# ---------------------------------------------------------------------------
function Set-StyleColor($styleName, $colorHex) {
    $s = $d.Styles($styleName)
    $r = [Convert]::ToInt32($colorHex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($colorHex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($colorHex.Substring(4,2), 16)
    $val = $b * 65536 + $g * 256 + $r
    $s.Font.Color = $val
}

Set-StyleColor "Heading 3" "002060"
Set-StyleColor "Heading 4" "002060"
Set-StyleColor "Heading 3 Char" "002060"
Set-StyleColor "Heading 4 Char" "002060"
